$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Line data")

$ws.Range("F2").Value = 47.5
$ws.Range("F3").Value = 47.5
